$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

$ws.Range("A18").Value = "20.10.18"
$ws.Range("B18").Value = 0.66666666666666663
$ws.Range("C18").Value = 0.79166666666666663
$ws.Range("E18").Value = "-Cascasde-Classifier Options and refactoring`n-Timer"
$ws.Range("E18").WrapText = $true

$ws.Rows.Item(18).RowHeight = 30

$ws.Range("C19").Select()
